$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# --- Row 9 (set first among labeled text rows we need for shared-string order: A9) ---
$ws.Range("A9").Value = "Расстояние"
$ws.Range("D9:H9").Value = 299000000000
$ws.Range("N9:R9").Value = 299000000000
$ws.Range("X9:AB9").Value = 299000000000
$ws.Range("AH9:AL9").Value = 299000000000

# --- Row 8 (values needed before row 10 formula) ---
$ws.Range("A8").Value = "Скорость"
$ws.Range("D8").Value = 15500
$ws.Range("E8").Value = 16000
$ws.Range("F8").Value = 16474.082208590102
$ws.Range("G8").Value = 16800
$ws.Range("H8").Value = 17000
$ws.Range("N8").Value = 15500
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 16474.082208590102
$ws.Range("Q8").Value = 16800
$ws.Range("R8").Value = 17000
$ws.Range("X8").Value = 15500
$ws.Range("Y8").Value = 16000
$ws.Range("Z8").Value = 16474.082208590102
$ws.Range("AA8").Value = 16800
$ws.Range("AB8").Value = 17000
$ws.Range("AH8").Value = 15500
$ws.Range("AI8").Value = 16000
$ws.Range("AJ8").Value = 16474.082208590102
$ws.Range("AK8").Value = 16800
$ws.Range("AL8").Value = 17000

# --- Row 10: A10 label then shared formula across D10:AN10, later cleared to sparse blocks ---
$ws.Range("A10").Value = "Потенциал"
$ws.Range("D10:AN10").Formula = "=D8*D8/D9"
$ws.Range("I10:M10").ClearContents()
$ws.Range("S10:W10").ClearContents()
$ws.Range("AC10:AG10").ClearContents()
$ws.Range("AM10:AN10").ClearContents()

# --- Row 11: A11 label, B11 styled-empty cell, then values per block ---
$ws.Range("A11").Value = "Энергия"
$ws.Range("B11").NumberFormat = "0.00E+00"

$ws.Range("D11").Value = [double]"-7.1334267735635203E+33"
$ws.Range("E11").Value = [double]"-7.1061679569910799E+33"
$ws.Range("F11").Value = [double]"-7.0792242163803102E+33"
$ws.Range("G11").Value = [double]"-7.0602038751635103E+33"
$ws.Range("H11").Value = [double]"-7.0484131362986404E+33"
$ws.Range("D11:H11").NumberFormat = "0.00E+00"

$ws.Range("N11").Value = [double]"-7.1269500046638803E+33"
$ws.Range("O11").Value = [double]"-7.0966565091819603E+33"
$ws.Range("P11").Value = [double]"-7.0693205980280099E+33"
$ws.Range("Q11").Value = [double]"-7.0504204913389295E+33"
$ws.Range("R11").Value = [double]"-7.03820959641974E+33"
$ws.Range("N11:R11").NumberFormat = "0.00E+00"

$ws.Range("X11").Value = [double]"-7.1835822843444601E+33"
$ws.Range("Y11").Value = [double]"-7.1549296200628795E+33"
$ws.Range("Z11").Value = [double]"-7.1285509052817498E+33"
$ws.Range("AA11").Value = [double]"-7.1089587593986302E+33"
$ws.Range("AB11").Value = [double]"-7.0970732281849996E+33"
$ws.Range("X11:AB11").NumberFormat = "0.00E+00"

$ws.Range("AH11").Value = [double]"-7.1152035793229804E+33"
$ws.Range("AI11").Value = [double]"-7.0908541707421002E+33"
$ws.Range("AJ11").Value = [double]"-7.0634764043825602E+33"
$ws.Range("AK11").Value = [double]"-7.0438200700252797E+33"
$ws.Range("AL11").Value = [double]"-7.0324463088538204E+33"
$ws.Range("AH11:AL11").NumberFormat = "0.00E+00"

# --- Row 13 text labels (F13, P13, Z13, AJ13) ---
$ws.Range("F13").Value = "Эйлера"
$ws.Range("P13").Value = "Эйлера-Крамера"
$ws.Range("Z13").Value = "Верле"
$ws.Range("AJ13").Value = "Бимана"

# --- Row 7: A7 label, then shared values per block (style 0.00E+00) ---
$ws.Range("A7").Value = "Начальная энергия"
$ws.Range("D7").Value = [double]"-7.1338311854762195E+33"
$ws.Range("E7").Value = [double]"-7.10626868547622E+33"
$ws.Range("F7").Value = [double]"-7.0787547374762195E+33"
$ws.Range("G7").Value = [double]"-7.0603486854762205E+33"
$ws.Range("H7").Value = [double]"-7.0485186854762203E+33"
$ws.Range("D7:H7").NumberFormat = "0.00E+00"

$ws.Range("N7").Value = [double]"-7.1338311854762195E+33"
$ws.Range("O7").Value = [double]"-7.10626868547622E+33"
$ws.Range("P7").Value = [double]"-7.0787547374762195E+33"
$ws.Range("Q7").Value = [double]"-7.0603486854762205E+33"
$ws.Range("R7").Value = [double]"-7.0485186854762203E+33"
$ws.Range("N7:R7").NumberFormat = "0.00E+00"

$ws.Range("X7").Value = [double]"-7.1338311854762195E+33"
$ws.Range("Y7").Value = [double]"-7.10626868547622E+33"
$ws.Range("Z7").Value = [double]"-7.0787547374762195E+33"
$ws.Range("AA7").Value = [double]"-7.0603486854762205E+33"
$ws.Range("AB7").Value = [double]"-7.0485186854762203E+33"
$ws.Range("X7:AB7").NumberFormat = "0.00E+00"

$ws.Range("AH7").Value = [double]"-7.1338311854762195E+33"
$ws.Range("AI7").Value = [double]"-7.10626868547622E+33"
$ws.Range("AJ7").Value = [double]"-7.0787547374762195E+33"
$ws.Range("AK7").Value = [double]"-7.0603486854762205E+33"
$ws.Range("AL7").Value = [double]"-7.0485186854762203E+33"
$ws.Range("AH7:AL7").NumberFormat = "0.00E+00"

# --- Row 12: A12 label, D12 standalone formula, E12:AL12 shared formula, then sparse-ify ---
$ws.Range("A12").Value = "Разность"
$ws.Range("D12").Formula = "=ABS(D7-D11)"
$ws.Range("E12:AL12").Formula = "=ABS(E7-E11)"
$ws.Range("I12:M12").ClearContents()
$ws.Range("S12:W12").ClearContents()
$ws.Range("AC12:AG12").ClearContents()

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 8.666666666666666
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666

# --- Selection ---
$ws.Range("G12").Select()
